# Edit pl_mw.xlsx values: update columns B, C, D, F, G, H, M, N for data rows 2-25
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @('B','C','D','F','G','H','M','N')

$data = @{
    "2" = @(0.6761305404375832,0.1628680422044511,0.007216698647265929,0.6759200295582559,0.5339648941648818,0.5969941988219176,0.8938952887598077,1.009474523385228)
    "3" = @(0.5937072249151925,0.1448294142447537,0.006983636417057681,0.6545518479094881,0.5118022218975113,0.5916943411649953,0.7953823004104521,1.026147084168585)
    "4" = @(0.5431234138413572,0.13368649741534,0.006839941306928665,0.6419647620210895,0.4986587300323322,0.5888388394540698,0.735426960661286,1.036889376269599)
    "5" = @(0.5225166324047166,0.1291288884313246,0.006781243687541405,0.6369687186559077,0.4934185108533455,0.5877751230767103,0.711122172333063,1.041393636395879)
    "6" = @(0.5190952981874659,0.1283710904278905,0.006771488753471999,0.6361471584869207,0.4925553468473396,0.5876045197617543,0.7070938999189025,1.042149206442909)
    "7" = @(0.5428454766141613,0.1336250997177615,0.006839150248383774,0.6418968449020497,0.498587590721371,0.5888240896275079,0.7350986697694424,1.036949609813131)
    "8" = @(0.6477061144960601,0.1566623097365039,0.007136467128251667,0.6684411841047222,0.5262262295796631,0.595083877738702,0.8598137416666276,1.01511813361113)
    "9" = @(0.8535249674827128,0.2013033141202811,0.007714447577555461,0.7247615522806115,0.5841566693833897,0.6105380840493524,1.108900719627187,0.9763315652394695)
    "10" = @(1.00486132638764,0.2337750088708788,0.008135597434797859,0.7687984681687539,0.629061108834378,0.6238543339753448,1.295099883073874,0.9503095634198058)
    "11" = @(1.073737015392624,0.2484766278136874,0.008326350915702108,0.7894219291131037,0.6500131542370582,0.6303436868291215,1.380598999299025,0.9390133595165153)
    "12" = @(1.099822912534478,0.2540336388264279,0.00839845773159098,0.7973173888935179,0.6580237425172299,0.6328635193489163,1.41309770270027,0.93481405051811)
    "13" = @(1.094204669336705,0.2528372914131012,0.008382934029953049,0.7956131308137628,0.6562950998590793,0.6323180446027834,1.406092975127137,0.9357149559564899)
    "14" = @(1.075883036784717,0.248934011016587,0.008332285785229487,0.7900697688917973,0.6506706509327103,0.630549741109661,1.383270197267507,0.9386663074721984)
    "15" = @(1.064661037837368,0.246541810719151,0.008301245447452033,0.7866854987420879,0.6472355075435416,0.6294747488772714,1.369306708120732,0.9404843090884611)
    "16" = @(1.000360745070509,0.2328128063974475,0.008123113854672681,0.7674626374550115,0.6277024921679981,0.6234389587694977,1.289528977173958,0.9510587303511873)
    "17" = @(0.9609224874187703,0.2243725145512201,0.00801361756575858,0.7558219599048073,0.6158547682277344,0.6198470366601327,1.240797156975077,0.9576847792816388)
    "18" = @(0.9382416843322972,0.2195113079946793,0.00795056056973209,0.7491821127603941,0.609089603960939,0.6178216713516633,1.212842379750754,0.9615468314796018)
    "19" = @(0.9305628895500035,0.2178642615361639,0.007929197471053584,0.7469434920818969,0.6068074770678891,0.6171428814981255,1.203389937073382,0.9628631879220251)
    "20" = @(0.9651204458577922,0.2252716792654326,0.008025281725828393,0.7570553733162768,0.6171108653690283,0.6202251968854853,1.245976980269148,0.9569741512005141)
    "21" = @(1.081264434704281,0.2500807761245767,0.008347165922472755,0.7916956532723987,0.6523206033184579,0.6310674366838782,1.389970426940437,0.9377972946244526)
    "22" = @(1.157195344241927,0.2662356040684415,0.008556789750908678,0.8148355974900738,0.6757786109755273,0.6385176819673291,1.484793481654194,0.9257208358184972)
    "23" = @(1.116667526426681,0.2576189381726124,0.008444980447940509,0.8024392969156509,0.6632174389174281,0.6345078954137762,1.434116660943033,0.9321243194713968)
    "24" = @(0.9632225713133948,0.2248651942070126,0.008020008688642832,0.7564975835463486,0.6165428398712436,0.6200541070856502,1.24363499024804,0.9572952622039921)
    "25" = @(0.7978247470900328,0.189284002399404,0.007558673130613869,0.7090627141532053,0.5680784676859787,0.6060145173194371,1.040987353227294,0.9863913484795663)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = "$($cols[$i])$row"
        $ws.Range($addr).Value = $values[$i]
    }
}
